$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (F1:H1) ---
$ws.Range("F1").Value = "FirstName"
$ws.Range("G1").Value = "Lastname"
$ws.Range("H1").Value = "EmployeeId"

# --- New data row 2 additions (F2:H2) ---
$ws.Range("F2").Value = "Katrina"
$ws.Range("G2").Value = "Kaif"
$ws.Range("H2").Value = 8677

# --- New data row 3 additions (F3:H3) ---
$ws.Range("F3").Value = "Kareena"
$ws.Range("G3").Value = "Kapoor"
$ws.Range("H3").Value = 8678

# --- New data row 3 (A3:D3) - set B/C/D before A so the shared-string
#     insertion order matches (Rahul Dravid, 12@misty@12, then Joele mathew) ---
$ws.Range("B3").Value = "Rahul Dravid"
$ws.Range("C3").Value = "12@misty@12"
$ws.Range("D3").Value = "12@misty@12"
$ws.Range("A3").Value = "Joele mathew`n"

# --- Selection moves to B8 ---
$ws.Range("B8").Select()
